$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the price series. It belongs chronologically
# between the existing row 33 and the former row 34, so insert a fresh row at
# position 34 (this shifts the old rows 34:43 down to 35:44, matching the diff).
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record's data.
$ws.Cells.Item(34, 1).Value = 3
$ws.Cells.Item(34, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).Value = 44876
$ws.Cells.Item(34, 5).Value = 5
$ws.Cells.Item(34, 6).Value = 300000000
$ws.Cells.Item(34, 7).Value = "Espárragos"
$ws.Cells.Item(34, 8).Value = "Verde"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 1090
$ws.Cells.Item(34, 11).Value = 1400
$ws.Cells.Item(34, 12).Value = 1500
$ws.Cells.Item(34, 13).Value = 1450
$ws.Cells.Item(34, 14).Value = "$/kilo"
$ws.Cells.Item(34, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(34, 16).Value = 1450
$ws.Cells.Item(34, 17).Value = 1
$ws.Cells.Item(34, 18).Value = "Hortaliza"
